$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$reqText = "LOM3202 -  Circuitos Elétricos  (Requisito)`n"
$labText = "LOM3221 -  Laboratório de Eletrônica  (Indicação de Conjunto)`n"

$ws.Range("B24").Value = $reqText
$ws.Range("C24").Value = $reqText

$ws.Range("B25").Value = $labText
$ws.Range("C25").Value = $labText
